$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 (position 1): "Somma Requests per DDI-first-ar" -> "Sum Requests for
# DDI-first-arch". Same physical worksheet, same row count (20 data rows),
# values refreshed for a new test run.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sum Requests for DDI-first-arch"

$ws1.Range("B1").Value = "Sum Requests for DDI-first-architecture-go3"
$ws1.Range("A3").Value = "Thu Sep 26 2024 01:03:21 GMT+0200 (Ora legale dell’Europa centrale)"
$ws1.Range("A4").Value = "Thu Sep 26 2024 01:13:11 GMT+0200 (Ora legale dell’Europa centrale)"
$ws1.Range("A5").Value = "Grain: Automatic"
$ws1.Range("A6").Value = "Aggregation type: Sum"
$ws1.Range("B11").Value = "DDI-first-architecture-go3, Requests (Sum), DDI-first-architecture-go3"

$ws1.Cells.Item(12,1).Value = 45561.044444444444
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(13,1).Value = 45561.04513888889
$ws1.Cells.Item(13,2).Value = 32
$ws1.Cells.Item(14,1).Value = 45561.04583333333
$ws1.Cells.Item(14,2).Value = 105
$ws1.Cells.Item(15,1).Value = 45561.04652777778
$ws1.Cells.Item(15,2).Value = 119
$ws1.Cells.Item(16,1).Value = 45561.04722222222
$ws1.Cells.Item(16,2).Value = 170
$ws1.Cells.Item(17,1).Value = 45561.04791666667
$ws1.Cells.Item(17,2).Value = 400
$ws1.Cells.Item(18,1).Value = 45561.04861111111
$ws1.Cells.Item(18,2).Value = 164
$ws1.Cells.Item(19,1).Value = 45561.049305555556
$ws1.Cells.Item(19,2).Value = 10
$ws1.Cells.Item(20,1).Value = 45561.05
$ws1.Cells.Item(20,2).Value = 0

# ---------------------------------------------------------------------------
# Sheet 2 (position 2): was "Media Response Time per DDI-fir" (18 data rows,
# Aggregation type: Media) -> becomes "Min Response Time for DDI-first" (20
# data rows, Aggregation type: Min). Two extra rows (19, 20) must appear.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Min Response Time for DDI-first"

$ws2.Range("B1").Value = "Min Response Time for DDI-first-architecture-go3"
$ws2.Range("A3").Value = "Thu Sep 26 2024 01:03:21 GMT+0200 (Ora legale dell’Europa centrale)"
$ws2.Range("A4").Value = "Thu Sep 26 2024 01:13:11 GMT+0200 (Ora legale dell’Europa centrale)"
$ws2.Range("A5").Value = "Grain: Automatic"
$ws2.Range("A6").Value = "Aggregation type: Min"
$ws2.Range("B11").Value = "DDI-first-architecture-go3, Response Time (Min), DDI-first-architecture-go3"

$ws2.Cells.Item(12,1).Value = 45561.044444444444
$ws2.Cells.Item(12,2).Value = 0
$ws2.Cells.Item(13,1).Value = 45561.04513888889
$ws2.Cells.Item(13,2).Value = 0.002
$ws2.Cells.Item(14,1).Value = 45561.04583333333
$ws2.Cells.Item(14,2).Value = 0.028
$ws2.Cells.Item(15,1).Value = 45561.04652777778
$ws2.Cells.Item(15,2).Value = 0.027
$ws2.Cells.Item(16,1).Value = 45561.04722222222
$ws2.Cells.Item(16,2).Value = 0.009
$ws2.Cells.Item(17,1).Value = 45561.04791666667
$ws2.Cells.Item(17,2).Value = 0.025
$ws2.Cells.Item(18,1).Value = 45561.04861111111
$ws2.Cells.Item(18,2).Value = 0.023

# Rows 19 & 20 did not exist on this physical sheet before (it only had 18
# rows) - add them, and make sure column A picks up the same date/time
# display format ("m/d/yy h:mm", Excel built-in numFmtId 22) as the rest of
# column A on this sheet so they share the same style index.
$ws2.Cells.Item(19,1).Value = 45561.049305555556
$ws2.Cells.Item(19,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(19,2).Value = 0.026
$ws2.Cells.Item(20,1).Value = 45561.05
$ws2.Cells.Item(20,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(20,2).Value = 0

# ---------------------------------------------------------------------------
# Sheet 3 (position 3): was "Min Response Time per DDI-first" (20 data rows,
# Aggregation type: Min) -> becomes "Avg Response Time for DDI-first" (18
# data rows, Aggregation type: Avg). Rows 19 & 20 are dropped entirely.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Avg Response Time for DDI-first"

$ws3.Range("B1").Value = "Avg Response Time for DDI-first-architecture-go3"
$ws3.Range("A3").Value = "Thu Sep 26 2024 01:03:21 GMT+0200 (Ora legale dell’Europa centrale)"
$ws3.Range("A4").Value = "Thu Sep 26 2024 01:13:11 GMT+0200 (Ora legale dell’Europa centrale)"
$ws3.Range("A5").Value = "Grain: Automatic"
$ws3.Range("A6").Value = "Aggregation type: Avg"
$ws3.Range("B11").Value = "DDI-first-architecture-go3, Response Time (Avg), DDI-first-architecture-go3"

$ws3.Cells.Item(12,1).Value = 45561.04513888889
$ws3.Cells.Item(12,2).Value = 0.1610625
$ws3.Cells.Item(13,1).Value = 45561.04583333333
$ws3.Cells.Item(13,2).Value = 0.050076190476190476
$ws3.Cells.Item(14,1).Value = 45561.04652777778
$ws3.Cells.Item(14,2).Value = 0.04663025210084034
$ws3.Cells.Item(15,1).Value = 45561.04722222222
$ws3.Cells.Item(15,2).Value = 0.06258235294117648
$ws3.Cells.Item(16,1).Value = 45561.04791666667
$ws3.Cells.Item(16,2).Value = 0.04234
$ws3.Cells.Item(17,1).Value = 45561.04861111111
$ws3.Cells.Item(17,2).Value = 0.041743902439024386
$ws3.Cells.Item(18,1).Value = 45561.049305555556
$ws3.Cells.Item(18,2).Value = 0.047

# Drop the now-unused rows 19 & 20 (this sheet previously had 20 rows).
$ws3.Range("A19:B20").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet 4 (position 4): "Max Response Time per DDI-first" -> "Max Response
# Time for DDI-first". Same physical worksheet, gains one new row (21).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Max Response Time for DDI-first"

$ws4.Range("B1").Value = "Max Response Time for DDI-first-architecture-go3"
$ws4.Range("A3").Value = "Thu Sep 26 2024 01:03:21 GMT+0200 (Ora legale dell’Europa centrale)"
$ws4.Range("A4").Value = "Thu Sep 26 2024 01:13:11 GMT+0200 (Ora legale dell’Europa centrale)"
$ws4.Range("A5").Value = "Grain: Automatic"
$ws4.Range("A6").Value = "Aggregation type: Max"
$ws4.Range("B11").Value = "DDI-first-architecture-go3, Response Time (Max), DDI-first-architecture-go3"

$ws4.Cells.Item(12,1).Value = 45561.044444444444
$ws4.Cells.Item(12,2).Value = 0
$ws4.Cells.Item(13,1).Value = 45561.04513888889
$ws4.Cells.Item(13,2).Value = 1.349
$ws4.Cells.Item(14,1).Value = 45561.04583333333
$ws4.Cells.Item(14,2).Value = 0.701
$ws4.Cells.Item(15,1).Value = 45561.04652777778
$ws4.Cells.Item(15,2).Value = 0.111
$ws4.Cells.Item(16,1).Value = 45561.04722222222
$ws4.Cells.Item(16,2).Value = 0.861
$ws4.Cells.Item(17,1).Value = 45561.04791666667
$ws4.Cells.Item(17,2).Value = 0.244
$ws4.Cells.Item(18,1).Value = 45561.04861111111
$ws4.Cells.Item(18,2).Value = 0.084
$ws4.Cells.Item(19,1).Value = 45561.049305555556
$ws4.Cells.Item(19,2).Value = 0.085
$ws4.Cells.Item(20,1).Value = 45561.05
$ws4.Cells.Item(20,2).Value = 0

# New row 21 - give column A the same date/time display as the rest of
# column A on this sheet (Excel built-in numFmtId 22) so it reuses the same
# style index instead of minting a new one.
$ws4.Cells.Item(21,1).Value = 45561.05069444444
$ws4.Cells.Item(21,1).NumberFormat = "m/d/yy h:mm"
$ws4.Cells.Item(21,2).Value = 0

# New column (21) picks up the same width as the other data columns on this
# sheet (11 characters).
$ws4.Columns.Item(21).ColumnWidth = 10.16666666667
